$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the whole data block B2:D9 to 0
$ws.Range("B2:D9").Value = 0

# Apply the two specific non-zero overrides from the diff
$ws.Range("C3").Value = -0.6872315258788352
$ws.Range("C7").Value = -0.6781206431808408
